$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts existing rows 25:69 down to 26:70)
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 44498
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = 100112032
$ws.Cells.Item(25, 7).Value = "Zapallo italiano"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 350
$ws.Cells.Item(25, 11).Value = 10000
$ws.Cells.Item(25, 12).Value = 11000
$ws.Cells.Item(25, 13).Value = 10571
$ws.Cells.Item(25, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(25, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(25, 16).Value = 176
$ws.Cells.Item(25, 17).Value = 60
$ws.Cells.Item(25, 18).Value = "Hortaliza"
